$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update greeting text for rule R10 (cell E8) from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Move active selection to E8
$ws.Range("E8").Select()
